# Fruta / hortaliza, semanal
#
# The underlying data rows (2-49) get their Fecha/Variedad/Volumen/Precio
# columns re-shuffled across rows (a row permutation): for each row we move
# the values of columns D (Fecha), H (Variedad), J (Volumen), K (Precio
# minimo), L (Precio maximo), M (Precio promedio ponderado) and P (Precio
# $/Kg) coming from a different "before" row into the "after" row, while all
# the other columns (Mercado, Region, Categoria, Calidad, Unidad de
# comercializacion, Origen, Kg o Unidades, Clasificacion, ...) stay put.
#
# Build the row-after -> row-before map first, snapshot every source value,
# then write all destinations - this two-phase approach is required because
# the mapping is a full permutation (every row is simultaneously a read
# source and a write destination).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowMap = @(
    @{After=2; Before=42},
    @{After=3; Before=23},
    @{After=4; Before=31},
    @{After=5; Before=46},
    @{After=6; Before=26},
    @{After=7; Before=30},
    @{After=8; Before=29},
    @{After=9; Before=37},
    @{After=10; Before=28},
    @{After=11; Before=7},
    @{After=12; Before=24},
    @{After=13; Before=43},
    @{After=14; Before=21},
    @{After=15; Before=25},
    @{After=16; Before=22},
    @{After=17; Before=5},
    @{After=18; Before=14},
    @{After=19; Before=32},
    @{After=20; Before=12},
    @{After=21; Before=20},
    @{After=22; Before=41},
    @{After=23; Before=13},
    @{After=24; Before=17},
    @{After=25; Before=27},
    @{After=26; Before=9},
    @{After=27; Before=36},
    @{After=28; Before=10},
    @{After=29; Before=18},
    @{After=30; Before=47},
    @{After=31; Before=34},
    @{After=32; Before=33},
    @{After=33; Before=45},
    @{After=34; Before=19},
    @{After=35; Before=38},
    @{After=36; Before=40},
    @{After=37; Before=44},
    @{After=38; Before=3},
    @{After=39; Before=8},
    @{After=40; Before=35},
    @{After=41; Before=11},
    @{After=42; Before=6},
    @{After=43; Before=4},
    @{After=44; Before=16},
    @{After=45; Before=49},
    @{After=46; Before=2},
    @{After=47; Before=15},
    @{After=48; Before=48},
    @{After=49; Before=39}
)

# Columns that travel with the permutation.
$cols = @(4, 8, 10, 11, 12, 13, 16)   # D, H, J, K, L, M, P

# Phase 1: snapshot the "before" values for every row that will act as a
# source, for every column that moves.
$snapshot = @{}
foreach ($entry in $rowMap) {
    $beforeRow = $entry.Before
    foreach ($col in $cols) {
        $key = "$beforeRow-$col"
        if (-not $snapshot.ContainsKey($key)) {
            $snapshot[$key] = $ws.Cells.Item($beforeRow, $col).Value2
        }
    }
}

# Phase 2: write the snapshotted values into their destination rows.
foreach ($entry in $rowMap) {
    $afterRow = $entry.After
    $beforeRow = $entry.Before
    foreach ($col in $cols) {
        $key = "$beforeRow-$col"
        $ws.Cells.Item($afterRow, $col).Value2 = $snapshot[$key]
    }
}
